$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# G8: "too many formulas" -> "didn't get clear on it"
$ws.Range("G8").Value = "didn't get clear on it"

# F9: new value "SVM"
$ws.Range("F9").Value = "SVM"

# G9: new value "didn't get clear on it" (same text as G8, shared string reused)
$ws.Range("G9").Value = "didn't get clear on it"

# Update the selection to G9
$ws.Range("G9").Select()
